# Applies the "archivo.xlsx" re-purposing edit:
#   - Sheet goes from an "ingreso/visita" log (columns A:K) to a
#     "examen" log (columns A:I) -- two trailing columns are removed.
#   - Header row text is replaced.
#   - The two data rows are replaced with new COVID-test records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two now-unused trailing columns (J:K) so the sheet's used
# range becomes A1:I3, shifting nothing else.
$ws.Range("J1:K3").Delete()

# --- Header row -----------------------------------------------------
$ws.Range("A1").Value = "ID Exámen"
$ws.Range("B1").Value = "Entidad de salud"
$ws.Range("C1").Value = "Tipo de documento"
$ws.Range("D1").Value = "Nro documento"
$ws.Range("E1").Value = "Nombres"
$ws.Range("F1").Value = "Apellidos"
$ws.Range("G1").Value = "Resultado"
$ws.Range("H1").Value = "Fecha del exámen"
$ws.Range("I1").Value = "Días de cuarentena"

# --- Row 2 ------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "COMPAÑIA DE MEDICINA PREPAGADA COLSANITAS S A"
$ws.Range("C2").Value = "Cédula de ciudadanía"
# Force-text (leading apostrophe) so digit-only strings don't get
# silently re-typed as numbers/dates by Excel.
$ws.Range("D2").Value = "'1143878531"
$ws.Range("E2").Value = "Victor Manuel"
$ws.Range("F2").Value = "Toro Cedeño"
$ws.Range("G2").Value = "Negativo"
$ws.Range("H2").Value = "'2020-12-03"
$ws.Range("I2").Value = "'0"

# --- Row 3 ------------------------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "COMPAÑIA DE MEDICINA PREPAGADA COLSANITAS S A"
$ws.Range("C3").Value = "Cédula de ciudadanía"
$ws.Range("D3").Value = "'1193474912"
$ws.Range("E3").Value = "Isabela"
$ws.Range("F3").Value = "Acevedo García"
$ws.Range("G3").Value = "Positivo"
$ws.Range("H3").Value = "'2020-12-02"
$ws.Range("I3").Value = "'15"
